$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-05 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-06 Wednesday", 2)
$d.Content.Find.Execute("476×2=952", $true, $false, $false, $false, $false, $true, 1, $false, "682×6=4092", 2)
$d.Content.Find.Execute("972×6=5832", $true, $false, $false, $false, $false, $true, 1, $false, "763×5=3815", 2)
$d.Content.Find.Execute("606×8=4848", $true, $false, $false, $false, $false, $true, 1, $false, "158×7=1106", 2)
$d.Content.Find.Execute("333×6=1998", $true, $false, $false, $false, $false, $true, 1, $false, "878×2=1756", 2)
$d.Content.Find.Execute("211×8=1688", $true, $false, $false, $false, $false, $true, 1, $false, "782×9=7038", 2)
$d.Content.Find.Execute("950×3=2850", $true, $false, $false, $false, $false, $true, 1, $false, "573×4=2292", 2)
$d.Content.Find.Execute("585×5=2925", $true, $false, $false, $false, $false, $true, 1, $false, "169×4=676", 2)
$d.Content.Find.Execute("911×4=3644", $true, $false, $false, $false, $false, $true, 1, $false, "707×2=1414", 2)
$d.Content.Find.Execute("785×9=7065", $true, $false, $false, $false, $false, $true, 1, $false, "979×6=5874", 2)
$d.Content.Find.Execute("763×7=5341", $true, $false, $false, $false, $false, $true, 1, $false, "258×3=774", 2)
$d.Content.Find.Execute("292×2=584", $true, $false, $false, $false, $false, $true, 1, $false, "933×5=4665", 2)
$d.Content.Find.Execute("416×8=3328", $true, $false, $false, $false, $false, $true, 1, $false, "625×7=4375", 2)
$d.Content.Find.Execute("702×3=2106", $true, $false, $false, $false, $false, $true, 1, $false, "490×2=980", 2)
$d.Content.Find.Execute("649×8=5192", $true, $false, $false, $false, $false, $true, 1, $false, "703×3=2109", 2)
$d.Content.Find.Execute("240×7=1680", $true, $false, $false, $false, $false, $true, 1, $false, "577×2=1154", 2)
$d.Content.Find.Execute("332×5=1660", $true, $false, $false, $false, $false, $true, 1, $false, "259×8=2072", 2)
$d.Content.Find.Execute("602×8=4816", $true, $false, $false, $false, $false, $true, 1, $false, "860×9=7740", 2)
$d.Content.Find.Execute("866×2=1732", $true, $false, $false, $false, $false, $true, 1, $false, "905×4=3620", 2)
$d.Content.Find.Execute("111×2=222", $true, $false, $false, $false, $false, $true, 1, $false, "679×8=5432", 2)
$d.Content.Find.Execute("165×5=825", $true, $false, $false, $false, $false, $true, 1, $false, "856×4=3424", 2)
$d.Content.Find.Execute("427×8=3416", $true, $false, $false, $false, $false, $true, 1, $false, "646×4=2584", 2)
$d.Content.Find.Execute("971×7=6797", $true, $false, $false, $false, $false, $true, 1, $false, "818×5=4090", 2)
$d.Content.Find.Execute("909×3=2727", $true, $false, $false, $false, $false, $true, 1, $false, "713×3=2139", 2)
$d.Content.Find.Execute("383×4=1532", $true, $false, $false, $false, $false, $true, 1, $false, "347×5=1735", 2)
$d.Content.Find.Execute("889×4=3556", $true, $false, $false, $false, $false, $true, 1, $false, "457×3=1371", 2)
